# "Section 1" sheet: add the missing Linkedin-Presence count for row 4,
# and update the saved cursor/selection position (view scrolled right to
# column C, active cell now K8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data point that was missing from the sample data (row 4, column B).
$ws.Range("B4").Value = 456

# Update the sheet's remembered selection/active cell to K8 (also scrolls
# the view so column C becomes the left-most visible column).
$ws.Range("K8").Select()
